$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 129 (ALC)
$ws.Range("H129").Value = 16638.61
$ws.Range("I129").Value = 556.5
$ws.Range("J129").Value = 21999.312
$ws.Range("K129").Value = 1669.5
$ws.Range("L129").Value = 65997.936
$ws.Range("M129").Value = 3330.5
$ws.Range("N129").Value = -75997.936

# Row 137 (ALC)
$ws.Range("H137").Value = 1415.6562
$ws.Range("I137").Value = 1099.56
$ws.Range("J137").Value = 2544.5715
$ws.Range("K137").Value = 3298.68
$ws.Range("L137").Value = 7633.7145
$ws.Range("M137").Value = -748.6799999999998
$ws.Range("N137").Value = -12733.7145

# Row 141 (ALC)
$ws.Range("H141").Value = 2151.568
$ws.Range("I141").Value = 1376.96
$ws.Range("J141").Value = 3170.7896
$ws.Range("K141").Value = 4130.88
$ws.Range("L141").Value = 9512.3688
$ws.Range("M141").Value = 1049.12
$ws.Range("N141").Value = -19872.3688

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 17501.52
$ws.Range("I32").Value = 19116.895
$ws.Range("J32").Value = 5655.4443
$ws.Range("K32").Value = 19116.895
$ws.Range("L32").Value = 5655.4443
$ws.Range("M32").Value = -18829.895
$ws.Range("N32").Value = -6229.4443

# Row 61 (ARM)
$ws.Range("H61").Value = 1709.5151
$ws.Range("I61").Value = 1300.5
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1300.5
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -1088.5
$ws.Range("N61").Value = -4424

# Row 74 (ARM)
$ws.Range("H74").Value = 864.31915
$ws.Range("I74").Value = 769.4
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 769.4
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = 104.6
$ws.Range("N74").Value = -4748

# Row 77 (ARM)
$ws.Range("H77").Value = 864.31915
$ws.Range("I77").Value = 769.4
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 3847
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = 521
$ws.Range("N77").Value = -23736

# Row 136 (ARM)
$ws.Range("H136").Value = 1709.5151
$ws.Range("I136").Value = 1300.5
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 3901.5
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -1351.5
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
# Row 38 (BSM)
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# Row 99 (BSM)
$ws.Range("H99").Value = 1500
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2
$ws.Range("N99").ClearContents()

# Row 134 (BSM)
$ws.Range("H134").Value = 44626.75
$ws.Range("I134").Value = 65234.5
$ws.Range("J134").Value = 3411.25
$ws.Range("K134").Value = 195703.5
$ws.Range("L134").Value = 10233.75
$ws.Range("M134").Value = -193168.5
$ws.Range("N134").Value = -15303.75

$ws = $wb.Worksheets.Item("CRP")
# Row 20 (CRP)
$ws.Range("H20").Value = 47555.555
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 47555.555
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 47555.555
$ws.Range("N20").Value = -48027.555

# Row 30 (CRP)
$ws.Range("H30").Value = 47555.555
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 47555.555
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 47555.555
$ws.Range("N30").Value = -47737.555

# Row 76 (CRP)
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 5000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -4685

# Row 79 (CRP)
$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 5000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -3908

# Row 99 (CRP)
$ws.Range("H99").Value = 32686.242
$ws.Range("I99").Value = 64545.875
$ws.Range("J99").Value = 2700.7058
$ws.Range("K99").Value = 64545.875
$ws.Range("L99").Value = 2700.7058
$ws.Range("M99").Value = -63047.875
$ws.Range("N99").Value = -5696.7058

# Row 126 (CRP)
$ws.Range("H126").Value = 32686.242
$ws.Range("I126").Value = 64545.875
$ws.Range("J126").Value = 2700.7058
$ws.Range("K126").Value = 193637.625
$ws.Range("L126").Value = 8102.117400000001
$ws.Range("M126").Value = -191167.625
$ws.Range("N126").Value = -13042.1174

# Row 128 (CRP)
$ws.Range("H128").Value = 47555.555
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 47555.555
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 47555.555
$ws.Range("N128").Value = -57515.555

# Row 134 (CRP)
$ws.Range("H134").Value = 914.6923
$ws.Range("I134").Value = 872.44446
$ws.Range("J134").Value = 1186.2858
$ws.Range("K134").Value = 2617.33338
$ws.Range("L134").Value = 3558.8574
$ws.Range("M134").Value = -82.33338000000003
$ws.Range("N134").Value = -8628.857400000001

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (CUL)
$ws.Range("H4").Value = 33333516
$ws.Range("I4").Value = 33333516
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 100000548
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -100000436

# Row 122 (CUL)
$ws.Range("H122").Value = 1545.5385
$ws.Range("I122").Value = 1714.8572
$ws.Range("J122").Value = 1348
$ws.Range("K122").Value = 15433.7148
$ws.Range("L122").Value = 12132
$ws.Range("M122").Value = -12983.7148
$ws.Range("N122").Value = -17032

# Row 129 (CUL)
$ws.Range("H129").Value = 41667350
$ws.Range("I129").Value = 910
$ws.Range("J129").Value = 166666670
$ws.Range("K129").Value = 2730
$ws.Range("L129").Value = 500000010
$ws.Range("M129").Value = 2270
$ws.Range("N129").Value = -500010010

$ws = $wb.Worksheets.Item("GSM")
# Row 4 (GSM)
$ws.Range("H4").Value = 70004
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 70004
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 70004
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -70228

# Row 126 (GSM)
$ws.Range("H126").Value = 4240
$ws.Range("I126").Value = 5533.3335
$ws.Range("J126").Value = 2300
$ws.Range("K126").Value = 16600.0005
$ws.Range("L126").Value = 6900
$ws.Range("M126").Value = -14130.0005
$ws.Range("N126").Value = -11840

# Row 136 (GSM)
$ws.Range("H136").Value = 17345
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 17345
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 52035
$ws.Range("N136").Value = -57135

$ws = $wb.Worksheets.Item("LTW")
# Row 68 (LTW)
$ws.Range("H68").Value = 1312.3
$ws.Range("I68").Value = 1077.875
$ws.Range("J68").Value = 2250
$ws.Range("K68").Value = 1077.875
$ws.Range("L68").Value = 2250
$ws.Range("M68").Value = -328.875
$ws.Range("N68").Value = -3748

# Row 71 (LTW)
$ws.Range("H71").Value = 1312.3
$ws.Range("I71").Value = 1077.875
$ws.Range("J71").Value = 2250
$ws.Range("K71").Value = 5389.375
$ws.Range("L71").Value = 11250
$ws.Range("M71").Value = -1645.375
$ws.Range("N71").Value = -18738

# Row 132 (LTW)
$ws.Range("H132").Value = 7953.0625
$ws.Range("I132").Value = 10217.782
$ws.Range("J132").Value = 2165.4443
$ws.Range("K132").Value = 30653.346
$ws.Range("L132").Value = 6496.3329
$ws.Range("M132").Value = -28123.346
$ws.Range("N132").Value = -11556.3329

$ws = $wb.Worksheets.Item("WVR")
# Row 54 (WVR)
$ws.Range("H54").Value = 232454.55
$ws.Range("I54").Value = 2500000
$ws.Range("J54").Value = 5700
$ws.Range("K54").Value = 2500000
$ws.Range("L54").Value = 5700
$ws.Range("M54").Value = -2499480
$ws.Range("N54").Value = -6740

# Row 122 (WVR)
$ws.Range("H122").Value = 30745.314
$ws.Range("I122").Value = 40930.152
$ws.Range("J122").Value = 1322.4445
$ws.Range("K122").Value = 122790.456
$ws.Range("L122").Value = 3967.3335
$ws.Range("M122").Value = -120340.456
$ws.Range("N122").Value = -8867.333500000001

# Row 132 (WVR)
$ws.Range("H132").Value = 1074.907
$ws.Range("I132").Value = 1005.39026
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 3016.17078
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -486.1707799999999
$ws.Range("N132").Value = -12560

# Row 136 (WVR)
$ws.Range("H136").Value = 4819.4116
$ws.Range("I136").Value = 5530.185
$ws.Range("J136").Value = 2077.8572
$ws.Range("K136").Value = 16590.555
$ws.Range("L136").Value = 6233.571599999999
$ws.Range("M136").Value = -14040.555
$ws.Range("N136").Value = -11333.5716
